$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 104
$ws.Range("I4").Value = 104
$ws.Range("K4").Value = 104
$ws.Range("M4").Value = 10
$ws.Range("H9").Value = 112.2
$ws.Range("I9").Value = 100
$ws.Range("K9").Value = 100
$ws.Range("M9").Value = 69
$ws.Range("H38").Value = 755.3333
$ws.Range("I38").Value = 42.57143
$ws.Range("K38").Value = 127.71429
$ws.Range("M38").Value = 244.28571
$ws.Range("H40").Value = 1100
$ws.Range("I40").Value = 900
$ws.Range("K40").Value = 900
$ws.Range("M40").Value = -725
$ws.Range("H43").Value = 2600
$ws.Range("J43").Value = 2600
$ws.Range("L43").Value = 2600
$ws.Range("N43").Value = -2738
$ws.Range("H51").Value = 6491
$ws.Range("I51").Value = 6491
$ws.Range("K51").Value = 6491
$ws.Range("M51").Value = -6007
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H98").Value = 1777.1538
$ws.Range("I98").Value = 1418.6364
$ws.Range("J98").Value = 3749
$ws.Range("K98").Value = 1418.6364
$ws.Range("L98").Value = 3749
$ws.Range("M98").Value = 79.36359999999991
$ws.Range("N98").Value = -6745
$ws.Range("H116").Value = 3564.6667
$ws.Range("I116").Value = 3564.6667
$ws.Range("K116").Value = 3564.6667
$ws.Range("M116").Value = -122.6667000000002
$ws.Range("H122").Value = 1777.1538
$ws.Range("I122").Value = 1418.6364
$ws.Range("J122").Value = 3749
$ws.Range("K122").Value = 4255.9092
$ws.Range("L122").Value = 11247
$ws.Range("M122").Value = -1805.9092
$ws.Range("N122").Value = -16147
$ws.Range("H132").Value = 4423.857
$ws.Range("I132").Value = 2916.4614
$ws.Range("K132").Value = 8749.3842
$ws.Range("M132").Value = -6219.3842
$ws.Range("H137").Value = 3040.7273
$ws.Range("I137").Value = 2927.5
$ws.Range("J137").Value = 3176.6
$ws.Range("K137").Value = 8782.5
$ws.Range("L137").Value = 9529.799999999999
$ws.Range("M137").Value = -6232.5
$ws.Range("N137").Value = -14629.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3857.1428
$ws.Range("I63").Value = 2750
$ws.Range("K63").Value = 2750
$ws.Range("M63").Value = -2064
$ws.Range("H66").Value = 3857.1428
$ws.Range("I66").Value = 2750
$ws.Range("K66").Value = 13750
$ws.Range("M66").Value = -10318
$ws.Range("H88").Value = 3616.8
$ws.Range("I88").Value = 1119
$ws.Range("J88").Value = 5282
$ws.Range("K88").Value = 1119
$ws.Range("L88").Value = 5282
$ws.Range("M88").Value = -713
$ws.Range("N88").Value = -6094
$ws.Range("H91").Value = 3616.8
$ws.Range("I91").Value = 1119
$ws.Range("J91").Value = 5282
$ws.Range("K91").Value = 1119
$ws.Range("L91").Value = 5282
$ws.Range("M91").Value = 285
$ws.Range("N91").Value = -8090
$ws.Range("H92").Value = 29666.666
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H122").Value = 2116.923
$ws.Range("I122").Value = 2116.923
$ws.Range("K122").Value = 6350.768999999999
$ws.Range("M122").Value = -3900.768999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 12666.667
$ws.Range("J35").Value = 12666.667
$ws.Range("L35").Value = 12666.667
$ws.Range("N35").Value = -13286.667
$ws.Range("H64").Value = 1178.1333
$ws.Range("J64").Value = 1005
$ws.Range("L64").Value = 1005
$ws.Range("N64").Value = -1455
$ws.Range("H67").Value = 1178.1333
$ws.Range("J67").Value = 1005
$ws.Range("L67").Value = 1005
$ws.Range("N67").Value = -2565
$ws.Range("H82").Value = 24320.6
$ws.Range("I82").Value = 8641.200000000001
$ws.Range("K82").Value = 8641.200000000001
$ws.Range("M82").Value = -8258.200000000001
$ws.Range("H85").Value = 24320.6
$ws.Range("I85").Value = 8641.200000000001
$ws.Range("K85").Value = 8641.200000000001
$ws.Range("M85").Value = -7315.200000000001
$ws.Range("H99").Value = 1113
$ws.Range("I99").Value = 1113
$ws.Range("K99").Value = 1113
$ws.Range("M99").Value = 385

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 724
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 748
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 748
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -1448
$ws.Range("H58").Value = 7013.25
$ws.Range("I58").Value = 5237
$ws.Range("K58").Value = 5237
$ws.Range("M58").Value = -5034
$ws.Range("H136").Value = 7013.25
$ws.Range("I136").Value = 5237
$ws.Range("K136").Value = 15711
$ws.Range("M136").Value = -13161

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 216.14285
$ws.Range("I40").Value = 162.16667
$ws.Range("K40").Value = 648.66668
$ws.Range("M40").Value = -579.66668
$ws.Range("H68").Value = 1685.8334
$ws.Range("H71").Value = 1685.8334
$ws.Range("H113").Value = 1699.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1699.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5098.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9438.5
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8931909
$ws.Range("I122").Value = 9618363
$ws.Range("K122").Value = 28855089
$ws.Range("M122").Value = -28852639

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6027.7144
$ws.Range("I40").Value = 3048.5
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 3048.5
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -2912.5
$ws.Range("N40").Value = -10272
$ws.Range("H46").Value = 3376.923
$ws.Range("I46").Value = 2842.8572
$ws.Range("K46").Value = 2842.8572
$ws.Range("M46").Value = -2654.8572

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 69499.5
$ws.Range("J68").Value = 69499.5
$ws.Range("L68").Value = 69499.5
$ws.Range("N68").Value = -71121.5
$ws.Range("H71").Value = 69499.5
$ws.Range("J71").Value = 69499.5
$ws.Range("L71").Value = 208498.5
$ws.Range("N71").Value = -216610.5
$ws.Range("H107").Value = 475.44446
$ws.Range("I107").Value = 509.875
$ws.Range("K107").Value = 1529.625
$ws.Range("M107").Value = 390.375
$ws.Range("H132").Value = 2162.8572
$ws.Range("I132").Value = 1662.4546
$ws.Range("K132").Value = 4987.3638
$ws.Range("M132").Value = -2457.3638
